$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7 - shifts existing rows 7:30 down to 8:31,
# carrying their formatting/values along.
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,J,K,Q,R,T are constant for this product/market
# across every row, so reuse those values; D/L/M/N/O/P/S are the new
# observations for this week.
$ws.Range("A7").Value2 = 1
$ws.Range("B7").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value2 = "Arica y Parinacota"
$ws.Range("D7").Value2 = 44498
$ws.Range("D7").NumberFormat = $ws.Range("D8").NumberFormat
$ws.Range("E7").Value2 = 15
$ws.Range("F7").Value2 = "Fruta"
$ws.Range("G7").Value2 = 100108
$ws.Range("H7").Value2 = "Tropicales y subtropicales"
$ws.Range("I7").Value2 = 100108001
$ws.Range("J7").Value2 = "Guayaba"
$ws.Range("K7").Value2 = "Sin especificar"
$ws.Range("L7").Value2 = "Segunda"
$ws.Range("M7").Value2 = 100
$ws.Range("N7").Value2 = 1200
$ws.Range("O7").Value2 = 1300
$ws.Range("P7").Value2 = 1250
$ws.Range("Q7").Value2 = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R7").Value2 = "Región de Arica y Parinacota"
$ws.Range("S7").Value2 = 1250
$ws.Range("T7").Value2 = 1
